$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117, shifting existing rows 117:213 down to 118:214.
$ws.Rows("117").Insert()

# Populate the new row 117 with the new price record (same market/region/
# category metadata as its neighbours, new variety-grade price data).
$ws.Range("A117").Value = 8
$ws.Range("B117").Value = "Terminal La Palmera de La Serena"
$ws.Range("C117").Value = "Coquimbo"
$ws.Range("D117").Value = 44680
$ws.Range("E117").Value = 4
$ws.Range("F117").Value = 100112021
$ws.Range("G117").Value = "Ají"
$ws.Range("H117").Value = "Inferno"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 200
$ws.Range("K117").Value = 34000
$ws.Range("L117").Value = 35000
$ws.Range("M117").Value = 34500
$ws.Range("N117").Value = "$/caja 25 kilos"
$ws.Range("O117").Value = "Provincia de Limarí"
$ws.Range("P117").Value = 1380
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"
